# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker data table (rows 16-38) is re-sorted: previously it was grouped
# primarily by worker (N° Doc Trabajador) and secondarily by period (Periodo
# Mora); now it is grouped primarily by period and secondarily by worker.
# The underlying (worker, period) -> (Valor Mora) values are unchanged - only
# the row order / which row holds which record changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for columns C (N° Doc Trabajador), D (Nombre Trabajador),
# E (Periodo Mora) and F (Valor Mora) for rows 16 through 38.
# Columns B (Tipo Doc Trabajador = "CC") and G (Salario Basico = 781242)
# are identical for every row and are left untouched.
$rows = @(
    @{ Row = 16; Doc = "10967447";   Nombre = "FRANCISCO AMADOR GUERRA GUERRERO"; Periodo = "1806"; Valor = 31249 },
    @{ Row = 17; Doc = "1049564184"; Nombre = "WILBERTO TORRES PEREZ";            Periodo = "1806"; Valor = 31249 },
    @{ Row = 18; Doc = "7922859";    Nombre = "JESUS ANTONIO MENDOZA GUERRERO";   Periodo = "1806"; Valor = 31249 },
    @{ Row = 19; Doc = "91077212";   Nombre = "SAMUEL LEON SUAREZ";               Periodo = "1806"; Valor = 5208  },

    @{ Row = 20; Doc = "10967447";   Nombre = "FRANCISCO AMADOR GUERRA GUERRERO"; Periodo = "1807"; Valor = 31249 },
    @{ Row = 21; Doc = "1049564184"; Nombre = "WILBERTO TORRES PEREZ";            Periodo = "1807"; Valor = 31249 },
    @{ Row = 22; Doc = "7922859";    Nombre = "JESUS ANTONIO MENDOZA GUERRERO";   Periodo = "1807"; Valor = 31249 },
    @{ Row = 23; Doc = "91077212";   Nombre = "SAMUEL LEON SUAREZ";               Periodo = "1807"; Valor = 31249 },

    @{ Row = 24; Doc = "10967447";   Nombre = "FRANCISCO AMADOR GUERRA GUERRERO"; Periodo = "1809"; Valor = 31249 },
    @{ Row = 25; Doc = "1049564184"; Nombre = "WILBERTO TORRES PEREZ";            Periodo = "1809"; Valor = 31249 },
    @{ Row = 26; Doc = "91077212";   Nombre = "SAMUEL LEON SUAREZ";               Periodo = "1809"; Valor = 31249 },

    @{ Row = 27; Doc = "10967447";   Nombre = "FRANCISCO AMADOR GUERRA GUERRERO"; Periodo = "1810"; Valor = 31249 },
    @{ Row = 28; Doc = "1049564184"; Nombre = "WILBERTO TORRES PEREZ";            Periodo = "1810"; Valor = 31249 },
    @{ Row = 29; Doc = "91077212";   Nombre = "SAMUEL LEON SUAREZ";               Periodo = "1810"; Valor = 31249 },

    @{ Row = 30; Doc = "10967447";   Nombre = "FRANCISCO AMADOR GUERRA GUERRERO"; Periodo = "1811"; Valor = 31249 },
    @{ Row = 31; Doc = "1049564184"; Nombre = "WILBERTO TORRES PEREZ";            Periodo = "1811"; Valor = 31249 },
    @{ Row = 32; Doc = "91077212";   Nombre = "SAMUEL LEON SUAREZ";               Periodo = "1811"; Valor = 31249 },

    @{ Row = 33; Doc = "10967447";   Nombre = "FRANCISCO AMADOR GUERRA GUERRERO"; Periodo = "1812"; Valor = 31249 },
    @{ Row = 34; Doc = "1049564184"; Nombre = "WILBERTO TORRES PEREZ";            Periodo = "1812"; Valor = 31249 },
    @{ Row = 35; Doc = "91077212";   Nombre = "SAMUEL LEON SUAREZ";               Periodo = "1812"; Valor = 31249 },

    @{ Row = 36; Doc = "10967447";   Nombre = "FRANCISCO AMADOR GUERRA GUERRERO"; Periodo = "1902"; Valor = 28124 },
    @{ Row = 37; Doc = "1049564184"; Nombre = "WILBERTO TORRES PEREZ";            Periodo = "1902"; Valor = 28124 },
    @{ Row = 38; Doc = "91077212";   Nombre = "SAMUEL LEON SUAREZ";               Periodo = "1902"; Valor = 28124 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Nombre
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Valor
}
